$d = $word.ActiveDocument

# Locate the "10,000" that appears in "... dus over een afstand van 10,000en lichtjaren" ...
$rng = $d.Content
$found = $rng.Find.Execute("10,000", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Toggling a character property (Bold on, then back off) forces Word to
    # split the surrounding run at the Find-match boundaries, so the
    # replacement text below lands in its own run instead of being merged
    # back into the neighbouring text.
    $rng.Bold = 1
    $rng.Text = "tienduizend"
    $rng.Bold = 0
}
